$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("eparcel_template")

# Row 2: COLIN MONTAGUE (existing row, overwrite in place)
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "COLIN MONTAGUE"
$ws.Range("D2").Value = "0242 578418"
$ws.Range("F2").Value = "1B Ribbonwood Place"
$ws.Range("I2").Value = "ALBION PARK RAIL"
$ws.Range("J2").Value = "New South Wales"
$ws.Range("K2").Value = 2527
$ws.Range("L2").Value = "121842358460-1639674302002"
$ws.Range("N2").Value = "cfm153"

# Row 3: Laura Gannaway
$ws.Range("A3").Value = 0.2
$ws.Range("B3").Value = "Laura Gannaway"
$ws.Range("D3").Value = "08 99216832"
$ws.Range("F3").Value = "1 Dayana Drive"
$ws.Range("I3").Value = "Geraldton"
$ws.Range("J3").Value = "Western Australia"
$ws.Range("K3").Value = 6530
$ws.Range("L3").Value = "121842365891-1640040448002"
$ws.Range("N3").Value = "lauden14"

# Row 4: Nikolas Taufatofua
$ws.Range("A4").Value = 0.2
$ws.Range("B4").Value = "Nikolas Taufatofua"
$ws.Range("D4").Value = "04 24543354"
$ws.Range("F4").Value = "66 Harts road"
$ws.Range("I4").Value = "Indooroopilly"
$ws.Range("J4").Value = "Queensland"
$ws.Range("K4").Value = 4068
$ws.Range("L4").Value = "121842365891-1640084522002"
$ws.Range("N4").Value = "raoul3t"

# Row 5: Alison Wood
$ws.Range("A5").Value = 0.4
$ws.Range("B5").Value = "Alison Wood"
$ws.Range("D5").Value = "0249 346748"
$ws.Range("F5").Value = "91 Ferraby Dr"
$ws.Range("I5").Value = "Metford"
$ws.Range("J5").Value = "New South Wales"
$ws.Range("K5").Value = 2323
$ws.Range("L5").Value = "121921850608-1640163399002"
$ws.Range("N5").Value = "madgecod"

# recipient_address2 (column G) is left blank for every row, but still
# present as an (empty) cell in the row - mirror row 2's existing blank
# G cell onto the newly added rows 3-5.
$ws.Range("G2").Copy($ws.Range("G3:G5"))
